# The product/environment data sheet stored column A as an Excel date
# serial (custom-formatted as "YYYY-MM-DD HH:MM:SS"). The dataloader now
# expects a plain YYYYMMDD integer instead, with no special date
# formatting, so convert every data row's date cell in place and strip
# the custom number-format style back to the workbook default ("Normal"
# / General), matching how the header row (style index 1) and the other
# plain numeric columns are already styled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Row 1 is the header ("날짜" / "상품번호" / 1..15); data starts on row 2.
for ($i = $firstRow + 1; $i -le $lastRow; $i++) {
    $cell = $ws.Cells.Item($i, 1)
    $dateValue = $cell.Value()
    $yyyymmdd = [int]$dateValue.ToString("yyyyMMdd")
    $cell.Value = $yyyymmdd
    $cell.Style = "Normal"
}
